$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell updates derived from the coin-price refresh diff.
# Numeric-looking text values (e.g. "0.9995") need an explicit
# text NumberFormat, otherwise Excel auto-converts them to numbers;
# ClearFormats() afterwards restores the original (default) cell style
# while keeping the stored value as text.

$ws.Range("D2").Value = "30.386.32"
$ws.Range("E2").Value = "  +0.34%  "
$ws.Range("D3").Value = "1.880.91"
$ws.Range("E3").Value = "  +0.69%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9995"
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = "  -0.10%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "245.47"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +4.57%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9994"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -0.10%  "
$ws.Range("E7").Value = "  +0.84%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2895"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  +1.45%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "42.73"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  +2.68%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.06540"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  -0.47%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "21.23"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  -0.56%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07778"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  -0.51%  "
$ws.Range("D13").Value = "1.878.67"
$ws.Range("E13").Value = "  +0.54%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.7344"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  +5.80%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "95.80"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  -1.14%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "5.147"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  +1.23%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "276.56"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  +2.79%  "
$ws.Range("D18").Value = "30.371.29"
$ws.Range("E18").Value = "  +0.75%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "13.41"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  -2.58%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.000007590"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  -1.59%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.001"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  +0.04%  "
$ws.Range("D22").Value = "2.124.79"
$ws.Range("E22").Value = "  +1.05%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.9996"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  -0.05%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "5.268"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  +0.22%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "6.185"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +0.49%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "9.271"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  -3.47%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "164.99"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  -0.84%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "19.00"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  +0.64%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.940"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  +0.06%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.387"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  +1.75%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.09908"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  +0.19%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.520"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  +4.39%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.335"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  -0.40%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.067"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  +0.34%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.04772"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  +0.65%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.128"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  -0.19%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.7006"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  -0.28%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.716"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  -0.10%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01851"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  -1.20%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.756"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  -0.55%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "6.439"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  +1.84%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "70.29"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  -3.67%  "
$ws.Range("B43").Value = "RenderToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.921"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -1.41%  "
$ws.Range("B44").Value = "TrustWalletToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.8449"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  +1.10%  "
$ws.Range("B45").Value = "TheSandbox"
$ws.Range("C45").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.4168"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  +0.04%  "
$ws.Range("B46").Value = "PaxDollar"
$ws.Range("C46").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.9995"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  -0.09%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "102.51"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -0.35%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "9.443"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  +3.65%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "7.117"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  -0.06%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "930.82"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  -5.25%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "35.37"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  +2.40%  "
